# Update column G ("K") values on Sheet1 per regenerated save_data.
# The underlying source data changed from a "Strike#" style count to
# a true strikeout ("K") count, so the G column values for rows 2-33
# are rewritten with their newly-computed totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 2
    6  = 1
    7  = 1
    8  = 4
    9  = 1
    10 = 7
    11 = 6
    12 = 7
    13 = 6
    14 = 5
    15 = 1
    16 = 1
    17 = 8
    18 = 9
    19 = 7
    20 = 10
    21 = 8
    22 = 7
    23 = 3
    24 = 3
    25 = 6
    26 = 1
    27 = 6
    28 = 6
    29 = 3
    30 = 4
    31 = 3
    32 = 5
    33 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
